$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "kemasan" is a new column inserted between "berat" (D) and "merk" (E),
# pushing "merk" from E1 to F1.
$ws.Range("F1").Value = "merk"
$ws.Range("E1").Value = "kemasan"

# Match the author's final active-cell selection.
[void]$ws.Range("H4").Select()
